# Append the three new log rows (25-27) to Sheet1, replicating the
# get_price logging format used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value into a cell as literal text, even when the text
# looks like a number/date/time/bool that Excel would otherwise parse
# into a typed value. We stage the text in a scratch cell that has been
# explicitly formatted as Text ("@"), copy it, and paste-special just the
# values into the destination - this carries the "stored as text" quality
# without leaving the destination cell's own number format touched.
function Set-TextValue($sheet, $cell, $value) {
    $helper = $sheet.Cells.Item(1048500, 16384)
    $helper.NumberFormat = "@"
    $helper.Value = $value
    $helper.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
}

$rows = @(
    @("2024-10-09 11:22:44", "get_price", "https://example.com/product", "Error fetching price: 'NoneType' object is not subscriptable", "2024-10-09", "11:22:44"),
    @("2024-10-09 11:22:47", "get_price", "https://example.com/product", "100.00", "2024-10-09", "11:22:47"),
    @("2024-10-09 11:22:53", "get_price", "https://example.com/product", "Error fetching price: 'NoneType' object is not subscriptable", "2024-10-09", "11:22:53")
)

# Columns whose sample values in this sheet are pure numbers/dates and so
# need the text-forcing helper (Result column D holds prices like
# "100.00"; Entered Date column E holds ISO dates like "2024-10-09").
$forceTextCols = @(4, 5)

$startRow = 25
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $value = $rowData[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($forceTextCols -contains $c) {
            Set-TextValue $ws $cell $value
        } else {
            $cell.Value = $value
        }
    }
}
